# Update 2025-02-25 FlashScore odds sheet with refreshed odds/correct-score values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 3.2
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 1.9
$ws.Range("U2").Value = 2.5
$ws.Range("V2").Value = 1.53
$ws.Range("Y2").Value = 1.5
$ws.Range("Z2").Value = 2.5
$ws.Range("AC2").Value = 10
$ws.Range("G3").Value = 2.45
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 2.8
$ws.Range("J3").Value = 3.1
$ws.Range("L3").Value = 3.4
$ws.Range("Q3").Value = 1.92
$ws.Range("R3").Value = 1.98
$ws.Range("AB3").Value = 13
$ws.Range("AD3").Value = 26
$ws.Range("AG3").Value = 11
$ws.Range("AK3").Value = 201
$ws.Range("AL3").Value = 11
$ws.Range("AM3").Value = 17
$ws.Range("AN3").Value = 12
$ws.Range("AO3").Value = 34
$ws.Range("AQ3").Value = 34
$ws.Range("Q4").Value = 2.01
$ws.Range("R4").Value = 1.89
$ws.Range("L6").Value = 3.9
$ws.Range("P6").Value = 3.55
$ws.Range("AB6").Value = 10.25
$ws.Range("AD6").Value = 18
$ws.Range("AE6").Value = 14.5
$ws.Range("AL6").Value = 11.5
$ws.Range("AP6").Value = 29
$ws.Range("M7").Value = 1.11
$ws.Range("N7").Value = 6.5
$ws.Range("Q7").Value = 2.6
$ws.Range("R7").Value = 1.48
$ws.Range("AR7").Value = 1.98
$ws.Range("AS7").Value = 1.88
$ws.Range("G8").Value = 2.7
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 3.4
$ws.Range("AO8").Value = 26
$ws.Range("G9").Value = 3
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 2.25
$ws.Range("J9").Value = 3.6
$ws.Range("L9").Value = 3
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.93
$ws.Range("R9").Value = 1.93
$ws.Range("U9").Value = 3.4
$ws.Range("V9").Value = 1.33
$ws.Range("W9").Value = 1.4
$ws.Range("X9").Value = 2.75
$ws.Range("Y9").Value = 1.73
$ws.Range("Z9").Value = 2
$ws.Range("AA9").Value = 10
$ws.Range("AB9").Value = 15
$ws.Range("AC9").Value = 11
$ws.Range("AK9").Value = 201
$ws.Range("AN9").Value = 9.5
$ws.Range("AP9").Value = 19
$ws.Range("AQ9").Value = 26
$ws.Range("G10").Value = 1.6
$ws.Range("I10").Value = 5.5
$ws.Range("K10").Value = 2.1
$ws.Range("M10").Value = 1.07
$ws.Range("N10").Value = 8.5
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("Y10").Value = 2.05
$ws.Range("Z10").Value = 1.7
$ws.Range("AF10").Value = 34
$ws.Range("AG10").Value = 8.5
$ws.Range("AI10").Value = 21
$ws.Range("AK10").Value = 501
$ws.Range("AM10").Value = 29
$ws.Range("AN10").Value = 19
$ws.Range("AO10").Value = 67
$ws.Range("AP10").Value = 51
$ws.Range("I11").Value = 2
$ws.Range("Y11").Value = 1.75
$ws.Range("Z11").Value = 2
$ws.Range("AA11").Value = 11
$ws.Range("AH11").Value = 7
$ws.Range("AK11").Value = 201
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("G13").Value = 1.75
$ws.Range("I13").Value = 3.9
$ws.Range("Y13").Value = 1.62
$ws.Range("Z13").Value = 2.2
$ws.Range("G14").Value = 1.42
$ws.Range("AL14").Value = 19
$ws.Range("AQ14").Value = 41
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = 1.65
$ws.Range("J15").Value = 4.5
$ws.Range("K15").Value = 2.3
$ws.Range("M15").Value = 1.01
$ws.Range("N15").Value = 15
$ws.Range("Q15").Value = 1.65
$ws.Range("R15").Value = 2.2
$ws.Range("W15").Value = 1.3
$ws.Range("X15").Value = 3.4
$ws.Range("Y15").Value = 1.67
$ws.Range("Z15").Value = 2.1
$ws.Range("AB15").Value = 26
$ws.Range("AF15").Value = 34
$ws.Range("AG15").Value = 15
$ws.Range("AJ15").Value = 41
$ws.Range("AL15").Value = 8.5
$ws.Range("AM15").Value = 9
$ws.Range("AK16").Value = 600
$ws.Range("N19").Value = 26
$ws.Range("Q21").Value = 2.1
$ws.Range("R21").Value = 1.7
$ws.Range("U21").Value = 3.75
$ws.Range("V21").Value = 1.25
$ws.Range("M22").Value = 1.05
$ws.Range("N22").Value = 8.5
$ws.Range("G23").Value = 3.75
$ws.Range("H23").Value = 3.65
$ws.Range("I23").Value = 1.88
$ws.Range("J23").Value = 4.05
$ws.Range("K23").Value = 2.2
$ws.Range("L23").Value = 2.45
$ws.Range("O23").Value = 1.26
$ws.Range("P23").Value = 3.65
$ws.Range("U23").Value = 2.82
$ws.Range("W23").Value = 1.38
$ws.Range("X23").Value = 2.92
$ws.Range("AC23").Value = 13
$ws.Range("AF23").Value = 37
$ws.Range("AH23").Value = 7.4
$ws.Range("AI23").Value = 14.5
$ws.Range("AK23").Value = 450
$ws.Range("AL23").Value = 7.7
$ws.Range("AO23").Value = 17.5
$ws.Range("AP23").Value = 15.5
